$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing accuracy values (rf_raw, rf_raw_with_imp_cols, rf_raw_one_hot_encoded)
$ws.Range("B2").Value = 0.68018433179723503
$ws.Range("B3").Value = 0.69216589861751154
$ws.Range("B4").Value = 0.69216589861751154

# Add new row for rf_cv approach
$ws.Range("A5").Value = "rf_cv"
$ws.Range("B5").Value = 0.68663594470046085

# Match the author's last selected cell
$ws.Range("D4").Select() | Out-Null
